$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (D) and Volume(1h) (E) columns store text exactly as scraped,
# matching the original inlineStr cell type (avoids Excel auto-converting
# numeric-looking strings like "507.05" or "1.00" into real numbers).
$ws.Columns.Item(4).NumberFormat = "@"
$ws.Columns.Item(5).NumberFormat = "@"

$ws.Range('D2').Value = '56.431.04'
$ws.Range('E2').Value = '  -1.48%  '
$ws.Range('D3').Value = '3.005.87'
$ws.Range('E3').Value = '  +0.33%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = '507.05'
$ws.Range('E5').Value = '  +0.97%  '
$ws.Range('D6').Value = '138.34'
$ws.Range('E6').Value = '  +0.45%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('E8').Value = '  -0.17%  '
$ws.Range('D9').Value = '7.11'
$ws.Range('E9').Value = '  -2.57%  '
$ws.Range('E10').Value = '  -0.37%  '
$ws.Range('D11').Value = '0.368'
$ws.Range('E11').Value = '  +2.95%  '
$ws.Range('D12').Value = '3.511.36'
$ws.Range('E12').Value = '  +0.39%  '
$ws.Range('E13').Value = '  -0.59%  '
$ws.Range('D14').Value = '25.33'
$ws.Range('E14').Value = '  -3.19%  '
$ws.Range('E15').Value = '  +1.57%  '
$ws.Range('D16').Value = '56.324.26'
$ws.Range('E16').Value = '  -1.67%  '
$ws.Range('D17').Value = '3.000.28'
$ws.Range('E17').Value = '  +0.34%  '
$ws.Range('E18').Value = '  -2.28%  '
$ws.Range('D19').Value = '12.92'
$ws.Range('E19').Value = '  +2.37%  '
$ws.Range('D20').Value = '8.00'
$ws.Range('E20').Value = '  +1.57%  '
$ws.Range('D21').Value = '332.68'
$ws.Range('E21').Value = '  +3.68%  '
$ws.Range('E22').Value = '  -0.21%  '
$ws.Range('D23').Value = '0.497'
$ws.Range('E23').Value = '  +0.75%  '
$ws.Range('D24').Value = '64.87'
$ws.Range('E24').Value = '  +2.57%  '
$ws.Range('D25').Value = '3.125.81'
$ws.Range('E25').Value = '  +0.54%  '
$ws.Range('B26').Value = 'Binance-PegBSC-USD'
$ws.Range('C26').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D26').Value = '1.00'
$ws.Range('E26').Value = '  +0.07%  '
$ws.Range('B27').Value = 'Kaspa'
$ws.Range('C27').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D27').Value = '0.165'
$ws.Range('E27').Value = '  +1.30%  '
$ws.Range('D28').Value = '0.0₃0943'
$ws.Range('E28').Value = '  +5.47%  '
$ws.Range('D29').Value = '6.35'
$ws.Range('E29').Value = '  -4.19%  '
$ws.Range('D30').Value = '6.85'
$ws.Range('E30').Value = '  -3.58%  '
$ws.Range('E31').Value = '  +0.49%  '
$ws.Range('D32').Value = '20.34'
$ws.Range('E32').Value = '  +0.62%  '
$ws.Range('E33').Value = '  -0.71%  '
$ws.Range('D34').Value = '152.57'
$ws.Range('E34').Value = '  -1.96%  '
$ws.Range('D35').Value = '4.43'
$ws.Range('E35').Value = '  -2.88%  '
$ws.Range('D36').Value = '5.82'
$ws.Range('E36').Value = '  +0.61%  '
$ws.Range('D37').Value = '26.34'
$ws.Range('E37').Value = '  +7.69%  '
$ws.Range('D38').Value = '1.22'
$ws.Range('E38').Value = '  -1.56%  '
$ws.Range('D39').Value = '0.0660'
$ws.Range('E39').Value = '  -0.60%  '
$ws.Range('D40').Value = '3.034.53'
$ws.Range('E40').Value = '  +0.36%  '
$ws.Range('D41').Value = '36.53'
$ws.Range('E41').Value = '  -3.49%  '
$ws.Range('E42').Value = '  +0.04%  '
$ws.Range('E43').Value = '  +1.06%  '
$ws.Range('E44').Value = '  +1.15%  '
$ws.Range('D45').Value = '2.198.17'
$ws.Range('E45').Value = '  +0.27%  '
$ws.Range('E46').Value = '  -2.87%  '
$ws.Range('B47').Value = 'VeChain'
$ws.Range('C47').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D47').Value = '0.0239'
$ws.Range('E47').Value = '  +1.93%  '
$ws.Range('B48').Value = 'ONDO'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D48').Value = '0.924'
$ws.Range('E48').Value = '  -1.28%  '
$ws.Range('D49').Value = '5.83'
$ws.Range('E49').Value = '  -2.12%  '
$ws.Range('D50').Value = '19.45'
$ws.Range('E50').Value = '  +0.79%  '
$ws.Range('D51').Value = '0.0850'
$ws.Range('E51').Value = '  -2.21%  '
